# "Updated symbol list" data refresh (scheduled scraper commit).
# Re-applies the latest Coin / Link / Price / Volume(1h) snapshot for every
# row whose values moved since the previous run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) cells store numeric-/percent-
# looking text (e.g. "244.82", "-0.87%"), exactly like the rest of the sheet.
# Force text formatting before writing them so Excel does not silently
# coerce them into real numbers/percentages.
$textCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7",
    "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12",
    "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17",
    "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22",
    "D23", "E23", "E24", "E25", "D26", "E26", "E27", "D28", "E28", "D40",
    "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45",
    "E46", "E47", "D48", "E48", "E49", "E50"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Per-cell value updates (Coin name, Link, Price, Volume(1h)).
$ws.Range("D2").Value = "244.82"
$ws.Range("E2").Value = "-0.87%"
$ws.Range("D3").Value = "27.49"
$ws.Range("E3").Value = "5.07%"
$ws.Range("E4").Value = "0.51%"
$ws.Range("D5").Value = "0.05685"
$ws.Range("E5").Value = "1.57%"
$ws.Range("D6").Value = "6.517"
$ws.Range("E6").Value = "0.58%"
$ws.Range("D7").Value = "0.8207"
$ws.Range("E7").Value = "0.74%"
$ws.Range("D8").Value = "0.8512"
$ws.Range("E8").Value = "0.88%"
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").Value = "0.06939"
$ws.Range("E9").Value = "-0.85%"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "0.02882"
$ws.Range("E10").Value = "2.27%"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "0.09392"
$ws.Range("E11").Value = "0.00%"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "0.001529"
$ws.Range("E12").Value = "0.41%"
$ws.Range("B13").Value = "CoinExToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D13").Value = "0.04021"
$ws.Range("E13").Value = "-13.55%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "0.0006021"
$ws.Range("E14").Value = "-93.93%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.006216"
$ws.Range("E15").Value = "1.40%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.514"
$ws.Range("E16").Value = "-2.66%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.012"
$ws.Range("E17").Value = "-0.24%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.319"
$ws.Range("E18").Value = "12.84%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3164"
$ws.Range("E19").Value = "1.65%"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").Value = "0.1332"
$ws.Range("E20").Value = "-0.35%"
$ws.Range("B21").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C21").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D21").Value = "0.03209"
$ws.Range("E21").Value = "0.35%"
$ws.Range("E22").Value = "-0.14%"
$ws.Range("D23").Value = "3.553"
$ws.Range("E23").Value = "-5.19%"
$ws.Range("E24").Value = "-0.12%"
$ws.Range("E25").Value = "-2.48%"
$ws.Range("D26").Value = "0.004478"
$ws.Range("E26").Value = "-2.00%"
$ws.Range("E27").Value = "22.93%"
$ws.Range("D28").Value = "0.0001407"
$ws.Range("E28").Value = "-27.51%"
$ws.Range("D40").Value = "0.03718"
$ws.Range("E40").Value = "1.53%"
$ws.Range("D41").Value = "0.005931"
$ws.Range("E41").Value = "73.70%"
$ws.Range("E42").Value = "-22.24%"
$ws.Range("D43").Value = "0.002351"
$ws.Range("E43").Value = "-10.17%"
$ws.Range("D44").Value = "0.009727"
$ws.Range("E44").Value = "21.13%"
$ws.Range("D45").Value = "0.00005101"
$ws.Range("E45").Value = "-5.20%"
$ws.Range("E46").Value = "-0.12%"
$ws.Range("E47").Value = "-30.43%"
$ws.Range("D48").Value = "0.002515"
$ws.Range("E48").Value = "3.53%"
$ws.Range("E49").Value = "-0.12%"
$ws.Range("E50").Value = "-0.12%"
